$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 0.001041666666666667
$ws.Range("K2").Value = 7268
$ws.Range("L2").Value = 0.014536
